$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 '24.678.20'
Set-TextCell $ws 2 5 '  -0.99%  '

Set-TextCell $ws 3 4 '1.658.95'
Set-TextCell $ws 3 5 '  -2.84%  '

Set-TextCell $ws 4 5 '  -0.32%  '

Set-TextCell $ws 5 4 '321.26'
Set-TextCell $ws 5 5 '  +2.68%  '

Set-TextCell $ws 6 4 '0.9996'
Set-TextCell $ws 6 5 '  +0.00%  '

Set-TextCell $ws 7 4 '0.3643'
Set-TextCell $ws 7 5 '  -2.65%  '

Set-TextCell $ws 8 4 '46.78'
Set-TextCell $ws 8 5 '  -5.42%  '

Set-TextCell $ws 9 4 '0.3267'
Set-TextCell $ws 9 5 '  -5.08%  '

Set-TextCell $ws 10 4 '1.134'
Set-TextCell $ws 10 5 '  -7.55%  '

Set-TextCell $ws 11 4 '0.07066'
Set-TextCell $ws 11 5 '  -6.29%  '

Set-TextCell $ws 12 4 '0.9987'
Set-TextCell $ws 12 5 '  -0.16%  '

Set-TextCell $ws 13 4 '5.995'
Set-TextCell $ws 13 5 '  -5.39%  '

Set-TextCell $ws 14 4 '19.56'
Set-TextCell $ws 14 5 '  -7.70%  '

Set-TextCell $ws 15 4 '1.658.44'
Set-TextCell $ws 15 5 '  -2.82%  '

Set-TextCell $ws 16 4 '6.626'
Set-TextCell $ws 16 5 '  -6.23%  '

Set-TextCell $ws 17 4 '0.00001047'
Set-TextCell $ws 17 5 '  -7.47%  '

Set-TextCell $ws 18 4 '0.06638'
Set-TextCell $ws 18 5 '  -1.40%  '

Set-TextCell $ws 19 4 '0.9988'
Set-TextCell $ws 19 5 '  +0.00%  '

Set-TextCell $ws 20 4 '78.91'
Set-TextCell $ws 20 5 '  -6.09%  '

Set-TextCell $ws 21 4 '5.933'
Set-TextCell $ws 21 5 '  -7.12%  '

Set-TextCell $ws 22 4 '15.76'
Set-TextCell $ws 22 5 '  -9.21%  '

Set-TextCell $ws 23 4 '12.65'
Set-TextCell $ws 23 5 '  -3.76%  '

Set-TextCell $ws 24 4 '24.646.72'
Set-TextCell $ws 24 5 '  -1.08%  '

Set-TextCell $ws 25 4 '2.472'
Set-TextCell $ws 25 5 '  +0.88%  '

Set-TextCell $ws 26 4 '2.398'
Set-TextCell $ws 26 5 '  -14.39%  '

Set-TextCell $ws 27 4 '148.01'
Set-TextCell $ws 27 5 '  -1.32%  '

Set-TextCell $ws 28 4 '18.63'
Set-TextCell $ws 28 5 '  -8.73%  '

Set-TextCell $ws 29 4 '1.841.97'
Set-TextCell $ws 29 5 '  -2.75%  '

Set-TextCell $ws 30 4 '1.215'
Set-TextCell $ws 30 5 '  -3.87%  '

Set-TextCell $ws 31 4 '125.32'
Set-TextCell $ws 31 5 '  -5.71%  '

Set-TextCell $ws 32 5 '  -3.62%  '

Set-TextCell $ws 33 4 '5.846'
Set-TextCell $ws 33 5 '  -14.25%  '

Set-TextCell $ws 34 4 '0.08479'
Set-TextCell $ws 34 5 '  -3.37%  '

Set-TextCell $ws 35 5 '  -5.22%  '

Set-TextCell $ws 36 4 '12.34'
Set-TextCell $ws 36 5 '  -10.74%  '

Set-TextCell $ws 37 4 '1.280'
Set-TextCell $ws 37 5 '  +0.64%  '

Set-TextCell $ws 38 4 '5.217'
Set-TextCell $ws 38 5 '  -7.36%  '

Set-TextCell $ws 41 4 '0.2079'
Set-TextCell $ws 41 5 '  -7.69%  '

Set-TextCell $ws 42 4 '8.216'
Set-TextCell $ws 42 5 '  -10.35%  '

Set-TextCell $ws 43 4 '0.9992'
Set-TextCell $ws 43 5 '  -0.04%  '

Set-TextCell $ws 44 4 '0.5939'
Set-TextCell $ws 44 5 '  -8.79%  '

Set-TextCell $ws 45 4 '3.860'
Set-TextCell $ws 45 5 '  +0.51%  '

Set-TextCell $ws 46 4 '12.75'
Set-TextCell $ws 46 5 '  -8.28%  '

Set-TextCell $ws 47 4 '0.5631'
Set-TextCell $ws 47 5 '  -8.85%  '

Set-TextCell $ws 48 4 '124.54'
Set-TextCell $ws 48 5 '  -3.69%  '

Set-TextCell $ws 49 4 '1.959'
Set-TextCell $ws 49 5 '  -7.72%  '

Set-TextCell $ws 50 4 '0.06968'
Set-TextCell $ws 50 5 '  -4.96%  '

Set-TextCell $ws 51 4 '1.197'
Set-TextCell $ws 51 5 '  -3.25%  '

# Row 39/40: coin identities swapped (VeChain now row 39, Hedera now row 40)
Set-TextCell $ws 39 2 'VeChain'
Set-TextCell $ws 39 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 39 4 '0.02238'
Set-TextCell $ws 39 5 '  -7.60%  '

Set-TextCell $ws 40 2 'Hedera'
Set-TextCell $ws 40 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 40 4 '0.06040'
Set-TextCell $ws 40 5 '  -9.41%  '
